$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values (formerly "Strike#"); update with newly
# regenerated values per the commit: "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals"
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 0
